# Add Pull Request Reference - Navya.Yogish
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$target = $ws.Range("B19")

# Insert the hyperlink pointing at the PR, using the URL as the
# display text for the OOXML <hyperlink display="..."> attribute.
$ws.Hyperlinks.Add(
    $target,
    "https://github.com/dhavalkeerthi/MRIInterns2026A/pull/11",
    "",
    "",
    "https://github.com/dhavalkeerthi/MRIInterns2026A/pull/11"
) | Out-Null

# Overwrite the cell text with the friendly PR title/description.
$target.Value = "demo by NavyaKKulal " + [char]0x00B7 + " Pull Request #11 " + [char]0x00B7 + " dhavalkeerthi/MRIInterns2026A"

# Match the selection left on the sheet after the edit.
$target.Select() | Out-Null
